$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 7 so the existing rows 7-12 shift down to 8-13,
# making room for a brand-new record at row 7.
$ws.Rows.Item(7).Insert()

# Row 7 becomes a new record - same market/product metadata, new date, volume, prices, origin.
$ws.Range("A7").Value = 8
$ws.Range("B7").Value = "Terminal La Palmera de La Serena"
$ws.Range("C7").Value = "Coquimbo"
$ws.Range("D7").Value = 44622
$ws.Range("D7").NumberFormat = $ws.Range("D8").NumberFormat
$ws.Range("E7").Value = 4
$ws.Range("F7").Value = "Fruta"
$ws.Range("G7").Value = 100104
$ws.Range("H7").Value = "Frutos de pepita"
$ws.Range("I7").Value = 100104003
$ws.Range("J7").Value = "Membrillo"
$ws.Range("K7").Value = "Champion"
$ws.Range("L7").Value = "Primera"
$ws.Range("M7").Value = 16
$ws.Range("N7").Value = 410000
$ws.Range("O7").Value = 420000
$ws.Range("P7").Value = 415000
$ws.Range("Q7").Value = "$/bins (450 kilos)"
$ws.Range("R7").Value = "Región de O'Higgins"
$ws.Range("S7").Value = 922
$ws.Range("T7").Value = 450
